# Applies the scheduled-runner price/profit refresh captured in the commit diff.
# Each row below corresponds to one Leve entry whose live Marketboard-derived
# columns (H/I/J/K/L/M/N) were recomputed; only the cells that actually changed
# value are touched, matching the unified diff exactly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 499
$ws.Range("I4").Value = 248.5
$ws.Range("K4").Value = 248.5
$ws.Range("M4").Value = -134.5

$ws.Range("H18").Value = 281.25
$ws.Range("I18").Value = 281.25
$ws.Range("K18").Value = 281.25
$ws.Range("M18").Value = 2.75

$ws.Range("H19").Value = 1796.2142
$ws.Range("I19").Value = 2202.5
$ws.Range("J19").Value = 1254.5
$ws.Range("K19").Value = 2202.5
$ws.Range("L19").Value = 1254.5
$ws.Range("M19").Value = -2027.5
$ws.Range("N19").Value = -1604.5

$ws.Range("H43").Value = 6662
$ws.Range("I43").Value = 6662
$ws.Range("K43").Value = 6662
$ws.Range("M43").Value = -6593

$ws.Range("H53").Value = 407.4
$ws.Range("J53").Value = 1000
$ws.Range("L53").Value = 1000
$ws.Range("N53").Value = -2274

$ws.Range("H92").Value = 455.77777
$ws.Range("J92").Value = 576.3333
$ws.Range("L92").Value = 576.3333
$ws.Range("N92").Value = -3072.3333

$ws.Range("H116").Value = 6997.75
$ws.Range("J116").Value = 7435.8
$ws.Range("L116").Value = 7435.8
$ws.Range("N116").Value = -14319.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1559.6666
$ws.Range("I2").Value = 1606.7858
$ws.Range("K2").Value = 1606.7858
$ws.Range("M2").Value = -1493.7858

$ws.Range("H32").Value = 5091.892
$ws.Range("I32").Value = 3284.879
$ws.Range("J32").Value = 19999.75
$ws.Range("K32").Value = 3284.879
$ws.Range("L32").Value = 19999.75
$ws.Range("M32").Value = -2997.879
$ws.Range("N32").Value = -20573.75

$ws.Range("H61").Value = 2283.7144
$ws.Range("I61").Value = 2164.6667
$ws.Range("J61").Value = 2998
$ws.Range("K61").Value = 2164.6667
$ws.Range("L61").Value = 2998
$ws.Range("M61").Value = -1952.6667
$ws.Range("N61").Value = -3422

$ws.Range("H116").Value = 1559.6666
$ws.Range("I116").Value = 1606.7858
$ws.Range("K116").Value = 1606.7858
$ws.Range("M116").Value = 687.2141999999999

$ws.Range("H122").Value = 5965.5557
$ws.Range("I122").Value = 6336.25
$ws.Range("K122").Value = 19008.75
$ws.Range("M122").Value = -16558.75

$ws.Range("H132").Value = 3471.8
$ws.Range("I132").Value = 2836.6667
$ws.Range("K132").Value = 8510.000100000001
$ws.Range("M132").Value = -5980.000100000001

$ws.Range("H136").Value = 2283.7144
$ws.Range("I136").Value = 2164.6667
$ws.Range("J136").Value = 2998
$ws.Range("K136").Value = 6494.000100000001
$ws.Range("L136").Value = 8994
$ws.Range("M136").Value = -3944.000100000001
$ws.Range("N136").Value = -14094

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1559.6666
$ws.Range("I3").Value = 1606.7858
$ws.Range("K3").Value = 1606.7858
$ws.Range("M3").Value = -1492.7858

$ws.Range("H64").Value = 2000
$ws.Range("J64").Value = 2000
$ws.Range("L64").Value = 2000
$ws.Range("N64").Value = -2450

$ws.Range("H67").Value = 2000
$ws.Range("J67").Value = 2000
$ws.Range("L67").Value = 2000
$ws.Range("N67").Value = -3560

$ws.Range("H80").Value = 1020.6667
$ws.Range("J80").Value = 876.25
$ws.Range("L80").Value = 876.25
$ws.Range("N80").Value = -2872.25

$ws.Range("H83").Value = 1020.6667
$ws.Range("J83").Value = 876.25
$ws.Range("L83").Value = 4381.25
$ws.Range("N83").Value = -14365.25

$ws.Range("H134").Value = 2135.762
$ws.Range("I134").Value = 1838.7858
$ws.Range("K134").Value = 5516.357400000001
$ws.Range("M134").Value = -2981.357400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2227.3333
$ws.Range("I16").Value = 2227.3333
$ws.Range("K16").Value = 2227.3333
$ws.Range("M16").Value = -1940.3333

$ws.Range("H107").Value = 2339.6365
$ws.Range("I107").Value = 1247
$ws.Range("J107").Value = 2582.4443
$ws.Range("K107").Value = 1247
$ws.Range("L107").Value = 2582.4443
$ws.Range("M107").Value = 673
$ws.Range("N107").Value = -6422.4443

$ws.Range("H113").Value = 2227.3333
$ws.Range("I113").Value = 2227.3333
$ws.Range("K113").Value = 2227.3333
$ws.Range("M113").Value = -57.33329999999978

$ws.Range("H122").Value = 3201.6924
$ws.Range("I122").Value = 2570.8572
$ws.Range("K122").Value = 7712.571599999999
$ws.Range("M122").Value = -5262.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 50000
$ws.Range("J51").Value = 50000
$ws.Range("L51").Value = 50000
$ws.Range("N51").Value = -51018

$ws.Range("H80").Value = 5885.778
$ws.Range("I80").Value = 5327.6665
$ws.Range("K80").Value = 5327.6665
$ws.Range("M80").Value = -4329.6665

$ws.Range("H83").Value = 5885.778
$ws.Range("I83").Value = 5327.6665
$ws.Range("K83").Value = 26638.3325
$ws.Range("M83").Value = -21646.3325

$ws.Range("H122").Value = 1799.8
$ws.Range("I122").Value = 1166.3334
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 3499.0002
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -1049.0002
$ws.Range("N122").Value = -13150

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4000
$ws.Range("I7").Value = 4000
$ws.Range("K7").Value = 4000
$ws.Range("M7").Value = -3888

$ws.Range("H122").Value = 3750.75
$ws.Range("J122").Value = 999
$ws.Range("L122").Value = 2997
$ws.Range("N122").Value = -7897

$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530

$ws.Range("H132").Value = 4047.3333
$ws.Range("I132").Value = 3998
$ws.Range("K132").Value = 11994
$ws.Range("M132").Value = -9464

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 17498.5
$ws.Range("I81").Value = 16000
$ws.Range("K81").Value = 32000
$ws.Range("M81").Value = -30939

$ws.Range("H84").Value = 17498.5
$ws.Range("I84").Value = 16000
$ws.Range("K84").Value = 160000
$ws.Range("M84").Value = -154696

$ws.Range("H107").Value = 800
$ws.Range("I107").Value = 800
$ws.Range("K107").Value = 2400
$ws.Range("M107").Value = -480

$ws.Range("H122").Value = 3705.125
$ws.Range("J122").Value = 4249.5
$ws.Range("L122").Value = 12748.5
$ws.Range("N122").Value = -17648.5

$ws.Range("H126").Value = 1999.6
$ws.Range("I126").Value = 1999.6
$ws.Range("K126").Value = 5998.799999999999
$ws.Range("M126").Value = -3528.799999999999

$ws.Range("H132").Value = 3826.8696
$ws.Range("I132").Value = 3325.4167
$ws.Range("J132").Value = 4373.909
$ws.Range("K132").Value = 9976.250100000001
$ws.Range("L132").Value = 13121.727
$ws.Range("M132").Value = -7446.250100000001
$ws.Range("N132").Value = -18181.727
